$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 43144
$ws.Range("C6").Value = 0.81388888888888899
$ws.Range("D6").Value = "Project Plan"
$ws.Range("E6").Value = "2700-Indu"
$ws.Range("F6").Value = "2100-Prerana"
$ws.Range("G6").Value = "Added the Task details and Dates"

$ws.Range("H9").Select()
